# Auto-generated Excel COM-interop edit script
# Applies the "Update automàtic: dades i banners [2026-02-06 18:34]" diff
# to the Dades_Meteo sheet: refreshed DATA_EXTRACCIO timestamps plus the
# small measurement deltas that came with that re-scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-06 18:33:21"
$ws.Range("O2").Value = "0.1 °C"
$ws.Range("E3").Value = "2026-02-06 18:33:23"
$ws.Range("E4").Value = "2026-02-06 18:33:26"
$ws.Range("J4").Value = "997.2 hPa"
$ws.Range("E5").Value = "2026-02-06 18:33:29"
$ws.Range("E6").Value = "2026-02-06 18:33:31"
$ws.Range("J6").Value = "998.6 hPa"
$ws.Range("E7").Value = "2026-02-06 18:33:33"
$ws.Range("H7").Value = "'60%"
$ws.Range("J7").Value = "998.3 hPa"
$ws.Range("E8").Value = "2026-02-06 18:33:36"
$ws.Range("E9").Value = "2026-02-06 18:33:38"
$ws.Range("E10").Value = "2026-02-06 18:33:40"
$ws.Range("H10").Value = "'88%"
$ws.Range("O10").Value = "9.4 °C"
$ws.Range("E11").Value = "2026-02-06 18:33:43"
$ws.Range("J11").Value = "998.6 hPa"
$ws.Range("E12").Value = "2026-02-06 18:33:45"
$ws.Range("H12").Value = "'58%"
$ws.Range("O12").Value = "14.3 °C"
$ws.Range("E13").Value = "2026-02-06 18:33:48"
$ws.Range("E14").Value = "2026-02-06 18:33:50"
$ws.Range("H14").Value = "'73%"
$ws.Range("E15").Value = "2026-02-06 18:33:52"
$ws.Range("J15").Value = "997.6 hPa"
$ws.Range("E16").Value = "2026-02-06 18:33:55"
$ws.Range("E17").Value = "2026-02-06 18:33:57"
$ws.Range("O17").Value = "6.2 °C"
$ws.Range("E18").Value = "2026-02-06 18:33:59"
$ws.Range("N18").Value = "-6.6 °C 17:59 TU"
$ws.Range("O18").Value = "-4.4 °C"
$ws.Range("E19").Value = "2026-02-06 18:34:01"
$ws.Range("I19").Value = "0.2 mm"
$ws.Range("J19").Value = "999.8 hPa"
$ws.Range("E20").Value = "2026-02-06 18:34:04"
$ws.Range("H20").Value = "'79%"
$ws.Range("E21").Value = "2026-02-06 18:34:06"
$ws.Range("J21").Value = "997.8 hPa"
$ws.Range("O21").Value = "8.8 °C"
$ws.Range("E22").Value = "2026-02-06 18:34:08"
$ws.Range("H22").Value = "'77%"
$ws.Range("E23").Value = "2026-02-06 18:34:11"
$ws.Range("H23").Value = "'82%"
$ws.Range("J23").Value = "997.6 hPa"
$ws.Range("E24").Value = "2026-02-06 18:34:13"
$ws.Range("J24").Value = "997.0 hPa"
$ws.Range("E25").Value = "2026-02-06 18:34:16"
$ws.Range("O25").Value = "4.6 °C"
$ws.Range("E26").Value = "2026-02-06 18:34:18"
$ws.Range("I26").Value = "0.4 mm"
$ws.Range("E27").Value = "2026-02-06 18:34:21"
$ws.Range("J27").Value = "997.6 hPa"
$ws.Range("O27").Value = "11.2 °C"
$ws.Range("E28").Value = "2026-02-06 18:34:23"
$ws.Range("H28").Value = "'82%"
$ws.Range("O28").Value = "5.1 °C"
$ws.Range("E29").Value = "2026-02-06 18:34:25"
$ws.Range("E30").Value = "2026-02-06 18:34:28"
$ws.Range("H30").Value = "'78%"
$ws.Range("E31").Value = "2026-02-06 18:34:30"
$ws.Range("J31").Value = "999.1 hPa"
$ws.Range("E32").Value = "2026-02-06 18:34:33"
$ws.Range("H32").Value = "'51%"
$ws.Range("J32").Value = "998.9 hPa"
$ws.Range("E33").Value = "2026-02-06 18:34:35"
$ws.Range("H33").Value = "'84%"
$ws.Range("E34").Value = "2026-02-06 18:34:37"
$ws.Range("H34").Value = "'73%"
$ws.Range("E35").Value = "2026-02-06 18:34:40"
$ws.Range("E36").Value = "2026-02-06 18:34:42"
$ws.Range("H36").Value = "'62%"
$ws.Range("J36").Value = "999.9 hPa"
$ws.Range("O36").Value = "13.1 °C"
